$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix regimen for POLARIS-2/3/4 rows: SOF/VEL -> SOF/VEL/VOX
$ws.Range("C15").Value = "SOF/VEL/VOX"
$ws.Range("C21").Value = "SOF/VEL/VOX"
$ws.Range("C22").Value = "SOF/VEL/VOX"

# Fix id typo and add display names for trial id rows
# (order matters for shared-string table insertion order)
$ws.Range("A25").Value = "GS-US-337-1468"
$ws.Range("B25").Value = "GS-US-337-1468 (LEPTON)"
$ws.Range("B24").Value = "GS-US-367-1871 (TRILOGY-3)"

# Widen display_name column
$ws.Columns("B").ColumnWidth = 28.83

# Update selection / scroll position
$ws.Range("C23").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1

# Match the author's window position (best effort)
$win.Left = 15180
$win.Top = 7760
